$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").Value = '''27.194.79'
$ws.Range("D2").Style = $s
$s = $ws.Range("E2").Style
$ws.Range("E2").Value = '''  +0.20%  '
$ws.Range("E2").Style = $s

$s = $ws.Range("D3").Style
$ws.Range("D3").Value = '''1.903.92'
$ws.Range("D3").Style = $s
$s = $ws.Range("E3").Style
$ws.Range("E3").Value = '''  +0.65%  '
$ws.Range("E3").Style = $s

$s = $ws.Range("D4").Style
$ws.Range("D4").Value = '''0.9999'
$ws.Range("D4").Style = $s
$s = $ws.Range("E4").Style
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = $s

$s = $ws.Range("D5").Style
$ws.Range("D5").Value = '''306.37'
$ws.Range("D5").Style = $s
$s = $ws.Range("E5").Style
$ws.Range("E5").Value = '''  -0.60%  '
$ws.Range("E5").Style = $s

$s = $ws.Range("D6").Style
$ws.Range("D6").Value = '''0.9997'
$ws.Range("D6").Style = $s
$s = $ws.Range("E6").Style
$ws.Range("E6").Value = '''  -0.03%  '
$ws.Range("E6").Style = $s

$s = $ws.Range("D7").Style
$ws.Range("D7").Value = '''0.5249'
$ws.Range("D7").Style = $s
$s = $ws.Range("E7").Style
$ws.Range("E7").Value = '''  +1.18%  '
$ws.Range("E7").Style = $s

$s = $ws.Range("D8").Style
$ws.Range("D8").Value = '''0.3775'
$ws.Range("D8").Style = $s
$s = $ws.Range("E8").Style
$ws.Range("E8").Value = '''  +1.40%  '
$ws.Range("E8").Style = $s

$s = $ws.Range("D9").Style
$ws.Range("D9").Value = '''0.07256'
$ws.Range("D9").Style = $s
$s = $ws.Range("E9").Style
$ws.Range("E9").Value = '''  +0.60%  '
$ws.Range("E9").Style = $s

$s = $ws.Range("D10").Style
$ws.Range("D10").Value = '''21.13'
$ws.Range("D10").Style = $s
$s = $ws.Range("E10").Style
$ws.Range("E10").Value = '''  +0.29%  '
$ws.Range("E10").Style = $s

$s = $ws.Range("D11").Style
$ws.Range("D11").Value = '''0.8992'
$ws.Range("D11").Style = $s
$s = $ws.Range("E11").Style
$ws.Range("E11").Value = '''  -0.69%  '
$ws.Range("E11").Style = $s

$s = $ws.Range("D12").Style
$ws.Range("D12").Value = '''0.08416'
$ws.Range("D12").Style = $s
$s = $ws.Range("E12").Style
$ws.Range("E12").Value = '''  +10.17%  '
$ws.Range("E12").Style = $s

$s = $ws.Range("D13").Style
$ws.Range("D13").Value = '''1.898.56'
$ws.Range("D13").Style = $s
$s = $ws.Range("E13").Style
$ws.Range("E13").Value = '''  +0.24%  '
$ws.Range("E13").Style = $s

$s = $ws.Range("D14").Style
$ws.Range("D14").Value = '''94.69'
$ws.Range("D14").Style = $s
$s = $ws.Range("E14").Style
$ws.Range("E14").Value = '''  -0.55%  '
$ws.Range("E14").Style = $s

$s = $ws.Range("E15").Style
$ws.Range("E15").Value = '''  -0.30%  '
$ws.Range("E15").Style = $s

$s = $ws.Range("D16").Style
$ws.Range("D16").Value = '''0.9998'
$ws.Range("D16").Style = $s
$s = $ws.Range("E16").Style
$ws.Range("E16").Value = '''  -0.04%  '
$ws.Range("E16").Style = $s

$s = $ws.Range("D17").Style
$ws.Range("D17").Value = '''0.000008605'
$ws.Range("D17").Style = $s
$s = $ws.Range("E17").Style
$ws.Range("E17").Value = '''  +1.10%  '
$ws.Range("E17").Style = $s

$s = $ws.Range("E18").Style
$ws.Range("E18").Value = '''  +1.21%  '
$ws.Range("E18").Style = $s

$s = $ws.Range("D19").Style
$ws.Range("D19").Value = '''0.9993'
$ws.Range("D19").Style = $s
$s = $ws.Range("E19").Style
$ws.Range("E19").Value = '''  -0.06%  '
$ws.Range("E19").Style = $s

$s = $ws.Range("D20").Style
$ws.Range("D20").Value = '''27.234.54'
$ws.Range("D20").Style = $s
$s = $ws.Range("E20").Style
$ws.Range("E20").Value = '''  +0.21%  '
$ws.Range("E20").Style = $s

$s = $ws.Range("D21").Style
$ws.Range("D21").Value = '''5.060'
$ws.Range("D21").Style = $s
$s = $ws.Range("E21").Style
$ws.Range("E21").Value = '''  +0.01%  '
$ws.Range("E21").Style = $s

$s = $ws.Range("D22").Style
$ws.Range("D22").Value = '''2.134.52'
$ws.Range("D22").Style = $s
$s = $ws.Range("E22").Style
$ws.Range("E22").Value = '''  -0.79%  '
$ws.Range("E22").Style = $s

$s = $ws.Range("D23").Style
$ws.Range("D23").Value = '''10.59'
$ws.Range("D23").Style = $s
$s = $ws.Range("E23").Style
$ws.Range("E23").Value = '''  -0.10%  '
$ws.Range("E23").Style = $s

$s = $ws.Range("D24").Style
$ws.Range("D24").Value = '''6.439'
$ws.Range("D24").Style = $s
$s = $ws.Range("E24").Style
$ws.Range("E24").Value = '''  -0.04%  '
$ws.Range("E24").Style = $s

$s = $ws.Range("D25").Style
$ws.Range("D25").Value = '''146.89'
$ws.Range("D25").Style = $s
$s = $ws.Range("E25").Style
$ws.Range("E25").Value = '''  +1.07%  '
$ws.Range("E25").Style = $s

$s = $ws.Range("E26").Style
$ws.Range("E26").Value = '''  +5.59%  '
$ws.Range("E26").Style = $s

$s = $ws.Range("D27").Style
$ws.Range("D27").Value = '''1.752'
$ws.Range("D27").Style = $s
$s = $ws.Range("E27").Style
$ws.Range("E27").Value = '''  -2.24%  '
$ws.Range("E27").Style = $s

$s = $ws.Range("E28").Style
$ws.Range("E28").Value = '''  +0.37%  '
$ws.Range("E28").Style = $s

$s = $ws.Range("D29").Style
$ws.Range("D29").Value = '''114.67'
$ws.Range("D29").Style = $s
$s = $ws.Range("E29").Style
$ws.Range("E29").Value = '''  +0.02%  '
$ws.Range("E29").Style = $s

$s = $ws.Range("D30").Style
$ws.Range("D30").Value = '''4.929'
$ws.Range("D30").Style = $s
$s = $ws.Range("E30").Style
$ws.Range("E30").Value = '''  -1.42%  '
$ws.Range("E30").Style = $s

$s = $ws.Range("D31").Style
$ws.Range("D31").Value = '''4.804'
$ws.Range("D31").Style = $s
$s = $ws.Range("E31").Style
$ws.Range("E31").Value = '''  -0.45%  '
$ws.Range("E31").Style = $s

$s = $ws.Range("D32").Style
$ws.Range("D32").Value = '''0.09284'
$ws.Range("D32").Style = $s
$s = $ws.Range("E32").Style
$ws.Range("E32").Value = '''  +0.65%  '
$ws.Range("E32").Style = $s

$s = $ws.Range("D33").Style
$ws.Range("D33").Value = '''0.8102'
$ws.Range("D33").Style = $s
$s = $ws.Range("E33").Style
$ws.Range("E33").Value = '''  +6.64%  '
$ws.Range("E33").Style = $s

$s = $ws.Range("D34").Style
$ws.Range("D34").Value = '''0.05064'
$ws.Range("D34").Style = $s
$s = $ws.Range("E34").Style
$ws.Range("E34").Value = '''  +0.01%  '
$ws.Range("E34").Style = $s

$s = $ws.Range("D35").Style
$ws.Range("D35").Value = '''1.236'
$ws.Range("D35").Style = $s
$s = $ws.Range("E35").Style
$ws.Range("E35").Value = '''  +3.24%  '
$ws.Range("E35").Style = $s

$s = $ws.Range("D36").Style
$ws.Range("D36").Value = '''2.953'
$ws.Range("D36").Style = $s
$s = $ws.Range("E36").Style
$ws.Range("E36").Value = '''  -2.42%  '
$ws.Range("E36").Style = $s

$s = $ws.Range("D37").Style
$ws.Range("D37").Value = '''3.356'
$ws.Range("D37").Style = $s
$s = $ws.Range("E37").Style
$ws.Range("E37").Value = '''  +2.39%  '
$ws.Range("E37").Style = $s

$s = $ws.Range("D38").Style
$ws.Range("D38").Value = '''2.610'
$ws.Range("D38").Style = $s
$s = $ws.Range("E38").Style
$ws.Range("E38").Value = '''  +1.86%  '
$ws.Range("E38").Style = $s

$s = $ws.Range("D39").Style
$ws.Range("D39").Value = '''0.5703'
$ws.Range("D39").Style = $s
$s = $ws.Range("E39").Style
$ws.Range("E39").Value = '''  +1.03%  '
$ws.Range("E39").Style = $s

$s = $ws.Range("D40").Style
$ws.Range("D40").Value = '''0.01991'
$ws.Range("D40").Style = $s
$s = $ws.Range("E40").Style
$ws.Range("E40").Value = '''  -0.37%  '
$ws.Range("E40").Style = $s

$s = $ws.Range("D41").Style
$ws.Range("D41").Value = '''1.072'
$ws.Range("D41").Style = $s
$s = $ws.Range("E41").Style
$ws.Range("E41").Value = '''  -0.65%  '
$ws.Range("E41").Style = $s

$s = $ws.Range("D42").Style
$ws.Range("D42").Value = '''6.655'
$ws.Range("D42").Style = $s
$s = $ws.Range("E42").Style
$ws.Range("E42").Value = '''  +0.82%  '
$ws.Range("E42").Style = $s

$s = $ws.Range("D43").Style
$ws.Range("D43").Value = '''8.969'
$ws.Range("D43").Style = $s
$s = $ws.Range("E43").Style
$ws.Range("E43").Value = '''  +0.54%  '
$ws.Range("E43").Style = $s

$s = $ws.Range("D44").Style
$ws.Range("D44").Value = '''117.78'
$ws.Range("D44").Style = $s

$s = $ws.Range("D45").Style
$ws.Range("D45").Value = '''0.1515'
$ws.Range("D45").Style = $s
$s = $ws.Range("E45").Style
$ws.Range("E45").Value = '''  +0.28%  '
$ws.Range("E45").Style = $s

$s = $ws.Range("D46").Style
$ws.Range("D46").Value = '''0.4842'
$ws.Range("D46").Style = $s
$s = $ws.Range("E46").Style
$ws.Range("E46").Value = '''  +0.37%  '
$ws.Range("E46").Style = $s

$s = $ws.Range("D47").Style
$ws.Range("D47").Value = '''0.9993'
$ws.Range("D47").Style = $s
$s = $ws.Range("E47").Style
$ws.Range("E47").Value = '''  -0.06%  '
$ws.Range("E47").Style = $s

$s = $ws.Range("E48").Style
$ws.Range("E48").Value = '''  -0.32%  '
$ws.Range("E48").Style = $s

$s = $ws.Range("D49").Style
$ws.Range("D49").Value = '''1.612'
$ws.Range("D49").Style = $s
$s = $ws.Range("E49").Style
$ws.Range("E49").Value = '''  +2.11%  '
$ws.Range("E49").Style = $s

$s = $ws.Range("D50").Style
$ws.Range("D50").Value = '''37.43'
$ws.Range("D50").Style = $s
$s = $ws.Range("E50").Style
$ws.Range("E50").Value = '''  +0.63%  '
$ws.Range("E50").Style = $s

$s = $ws.Range("D51").Style
$ws.Range("D51").Value = '''63.59'
$ws.Range("D51").Style = $s
$s = $ws.Range("E51").Style
$ws.Range("E51").Value = '''  -0.07%  '
$ws.Range("E51").Style = $s
